# Generate Report for Archive
#
# 1. Update the status text "Ready for handoff" -> "In Translation" everywhere
#    it is used (Overview sheet columns E/F, and column C on each language
#    sheet).
# 2. Shrink the "Status" column(s) that held that text, since the new text
#    is shorter (mirrors the auto-fit that produced the narrower column
#    width in the original workbook).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
